$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "H 72" record) so all subsequent rows shift up by one.
$ws.Rows(2).Delete()
